$wb = $excel.ActiveWorkbook

# --- DATA_FILE sheet (sheet1) ---
$ws1 = $wb.Worksheets.Item("DATA_FILE")

# Clear row 10's data (TC_009_CLICK_ON_HELP / Chrome) - leaves A10 with its
# existing style but no value, and removes B10 entirely (it had no style).
$ws1.Range("A10:B10").ClearContents()

# Remove the now-unused last row (row 19) completely.
$ws1.Rows("19:19").Delete()

# Update DATA_FILE's own selection to the (now-empty) row 10.
$ws1.Range("A10:XFD10").Select()

# --- API_DATA_FILE sheet (sheet2) becomes the active / selected sheet ---
$ws2 = $wb.Worksheets.Item("API_DATA_FILE")
$ws2.Activate()
$ws2.Range("B8").Select()
